$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 8.073476000000001
$ws.Range("H2").Value = 24.220428
$ws.Range("I2").Value = 0.2745703625077554
$ws.Range("J2").Value = 0.2745703625077553
$ws.Range("Q2").Value = 0.151886303988
$ws.Range("R2").Value = 1.366976735892
$ws.Range("S2").Value = 0.2745703625077554
$ws.Range("T2").Value = 0.2745703625077553

# Row 3
$ws.Range("I3").Value = 0.3498220011746915
$ws.Range("J3").Value = 0.3498220011746914
$ws.Range("S3").Value = 0.3498220011746915
$ws.Range("T3").Value = 0.3498220011746914

# Row 4
$ws.Range("G4").Value = 11.04437933333333
$ws.Range("H4").Value = 33.133138
$ws.Range("I4").Value = 0.3756076363175532
$ws.Range("J4").Value = 0.3756076363175532
$ws.Range("Q4").Value = 0.207777908398
$ws.Range("R4").Value = 1.870001175582
$ws.Range("S4").Value = 0.3756076363175532
$ws.Range("T4").Value = 0.3756076363175532
